$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.656.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.424.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.66%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.512"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.49%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  +8.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.96%  "

$ws.Range("E11").Value = "  +1.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.95%  "

$ws.Range("E13").Value = "  -2.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.801.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.424.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.833"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.504.10"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.29%  "

$ws.Range("E20").Value = "  +2.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0922"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.10%  "

$ws.Range("E25").Value = "  +1.97%  "

$ws.Range("E27").Value = "  +3.13%  "

$ws.Range("E28").Value = "  -3.58%  "

$ws.Range("E29").Value = "  +2.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.02%  "

$ws.Range("E32").Value = "  +21.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.73%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0780"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.10%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.72%  "

$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("E37").Value = "  +3.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.88%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "120.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.58%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.109"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.75%  "

$ws.Range("E44").Value = "  +4.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.941.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.68%  "

$ws.Range("E46").Value = "  +1.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.98%  "

$ws.Range("E48").Value = "  +2.99%  "

$ws.Range("E49").Value = "  +11.14%  "

$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.83%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.54%  "
